$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the LEG1_DIRECTION / LEG2_DIRECTION values on the data row.
$ws.Range("AL2").Value = "R"
$ws.Range("AM2").Value = "P"

# Update the view: scroll right and move the selection.
$ws.Application.ActiveWindow.ScrollColumn = 28
$ws.Range("AM2").Select()
